$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.792.54"
$ws.Range("E2").Value = "'  +1.60%  "
$ws.Range("D3").Value = "'1.708.55"
$ws.Range("E3").Value = "'  +1.77%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "'  +0.46%  "
$ws.Range("D5").Value = "'311.08"
$ws.Range("E5").Value = "'  +1.43%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "'  +0.25%  "
$ws.Range("D7").Value = "'0.3748"
$ws.Range("E7").Value = "'  +1.11%  "
$ws.Range("D8").Value = "'49.66"
$ws.Range("E8").Value = "'  +3.32%  "
$ws.Range("D9").Value = "'0.3438"
$ws.Range("E9").Value = "'  -0.04%  "
$ws.Range("D10").Value = "'1.205"
$ws.Range("E10").Value = "'  +2.04%  "
$ws.Range("D11").Value = "'0.07538"
$ws.Range("E11").Value = "'  +3.87%  "
$ws.Range("D12").Value = "'1.004"
$ws.Range("E12").Value = "'  +0.49%  "
$ws.Range("D13").Value = "'21.06"
$ws.Range("E13").Value = "'  +3.35%  "
$ws.Range("D14").Value = "'6.298"
$ws.Range("E14").Value = "'  +3.20%  "
$ws.Range("D15").Value = "'7.040"
$ws.Range("E15").Value = "'  +4.34%  "
$ws.Range("D16").Value = "'1.715.27"
$ws.Range("E16").Value = "'  +2.30%  "
$ws.Range("D17").Value = "'0.00001134"
$ws.Range("E17").Value = "'  +2.25%  "
$ws.Range("D18").Value = "'0.06729"
$ws.Range("E18").Value = "'  +0.12%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "'  +0.31%  "
$ws.Range("D20").Value = "'84.64"
$ws.Range("E20").Value = "'  +4.35%  "
$ws.Range("D21").Value = "'17.31"
$ws.Range("E21").Value = "'  +5.35%  "
$ws.Range("D22").Value = "'6.389"
$ws.Range("E22").Value = "'  +4.77%  "
$ws.Range("D23").Value = "'13.21"
$ws.Range("E23").Value = "'  +10.62%  "
$ws.Range("D24").Value = "'24.832.93"
$ws.Range("E24").Value = "'  +2.04%  "
$ws.Range("D25").Value = "'2.452"
$ws.Range("E25").Value = "'  +0.81%  "
$ws.Range("D26").Value = "'2.785"
$ws.Range("E26").Value = "'  +4.69%  "
$ws.Range("D27").Value = "'20.43"
$ws.Range("E27").Value = "'  +4.33%  "
$ws.Range("D28").Value = "'152.04"
$ws.Range("E28").Value = "'  -0.09%  "
$ws.Range("D29").Value = "'132.60"
$ws.Range("E29").Value = "'  +4.33%  "
$ws.Range("D30").Value = "'1.902.53"
$ws.Range("E30").Value = "'  +2.15%  "
$ws.Range("D31").Value = "'1.235"
$ws.Range("E31").Value = "'  +27.60%  "
$ws.Range("D32").Value = "'6.908"
$ws.Range("E32").Value = "'  +9.38%  "
$ws.Range("D33").Value = "'4.230"
$ws.Range("E33").Value = "'  +5.02%  "
$ws.Range("D34").Value = "'1.845"
$ws.Range("E34").Value = "'  +6.02%  "
$ws.Range("D35").Value = "'13.80"
$ws.Range("E35").Value = "'  +12.15%  "
$ws.Range("D36").Value = "'0.08793"
$ws.Range("E36").Value = "'  +3.70%  "
$ws.Range("D37").Value = "'5.608"
$ws.Range("E37").Value = "'  +5.11%  "
$ws.Range("D38").Value = "'0.06703"
$ws.Range("E38").Value = "'  +3.34%  "
$ws.Range("D39").Value = "'9.307"
$ws.Range("E39").Value = "'  +2.75%  "
$ws.Range("D40").Value = "'0.02411"
$ws.Range("E40").Value = "'  +3.31%  "
$ws.Range("D41").Value = "'0.2239"
$ws.Range("E41").Value = "'  +6.11%  "
$ws.Range("D42").Value = "'1.273"
$ws.Range("E42").Value = "'  +1.12%  "
$ws.Range("D43").Value = "'0.6447"
$ws.Range("E43").Value = "'  +4.44%  "
$ws.Range("D44").Value = "'0.9991"
$ws.Range("E44").Value = "'  +0.16%  "
$ws.Range("D45").Value = "'13.95"
$ws.Range("E45").Value = "'  +7.66%  "
$ws.Range("D46").Value = "'0.6162"
$ws.Range("E46").Value = "'  +3.65%  "
$ws.Range("D47").Value = "'3.827"
$ws.Range("E47").Value = "'  +1.36%  "
$ws.Range("D48").Value = "'2.138"
$ws.Range("E48").Value = "'  +5.63%  "
$ws.Range("D49").Value = "'130.03"
$ws.Range("E49").Value = "'  +2.24%  "
$ws.Range("D50").Value = "'0.07318"
$ws.Range("E50").Value = "'  +1.48%  "
$ws.Range("D51").Value = "'79.99"
$ws.Range("E51").Value = "'  +5.49%  "
